# Publish IG 1.0.1
# - Strip the "id: " label prefix from the Identifier value
# - Bump Version from 1.0.0 to 1.0.1
# - Replace the Contact value with the MedCom contact string
# - Insert a new "Jurisdiction" metadata row (with an empty value) right
#   after the "Contact" row, pushing the remaining rows down by one

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Identifier (row 3): "id: 1.2.208.184.100.10" -> "1.2.208.184.100.10"
$ws.Range("B3").Value = "1.2.208.184.100.10"

# Version (row 4): "1.0.0" -> "1.0.1"
$ws.Range("B4").Value = "1.0.1"

# Contact (row 11): "No display for ContactDetail" -> "MedCom (http://www.medcom.dk)"
$ws.Range("B11").Value = "MedCom (http://www.medcom.dk)"

# Insert a new row 12 for "Jurisdiction" with an empty value, matching the
# formatting of the surrounding data rows.
$ws.Rows.Item(12).Insert()
$ws.Range("A13:B13").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

$excel.CutCopyMode = 0
